$wb = $excel.ActiveWorkbook

$wsQ = $wb.Worksheets.Item("Questions")
$wsA = $wb.Worksheets.Item("Answers")

# Expand the Answers table so it covers the four new answer rows
$lo = $wsA.ListObjects.Item(1)
$lo.Resize($wsA.Range("A1:C257"))

# Fill in the missing answers for question 7 (matches the shared-string
# insertion order produced by the original edit: Guet, Faire pousser,
# Défilé, Marathon)
$wsA.Range("A256").Value = "Guet de la Cathédrale"
$wsA.Range("B256").Value = 1
$wsA.Range("C256").Value = 7

$wsA.Range("A255").Value = "Faire pousser un palmier"
$wsA.Range("B255").Value = 0
$wsA.Range("C255").Value = 7

$wsA.Range("A254").Value = "Défilé militaire"
$wsA.Range("B254").Value = 0
$wsA.Range("C254").Value = 7

$wsA.Range("A257").Value = "Marathon de Lausanne"
$wsA.Range("B257").Value = 0
$wsA.Range("C257").Value = 7

# Restore the saved scroll/selection state for both sheets
$wsQ.Activate()
$wsQ.Range("A9").Select()

$wsA.Activate()
$excel.Goto($wsA.Range("A241"), $true)
$wsA.Range("A261").Select()
